# Update for November 2025
# Applies the content changes described by the commit: refreshed manuscript /
# report / author / undergraduate counts, the updated journal-report table,
# the expanded undergraduate roster, and the Author/Date style cleanup.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- "Manuscripts and Reports" section ---------------------------------
Replace-Text "2389 Manuscripts and 4440 Reports" "2791 Manuscripts and 4470 Reports"
Replace-Text "In all, 4440 have been prepared" "In all, 4470 have been prepared"

# Table of reports-by-journal
Replace-Text "511" "572"
Replace-Text "917" "973"
Replace-Text "600" "609"
Replace-Text "454" "485"
Replace-Text "262" "258"
Replace-Text "943" "967"
Replace-Text "244" "260"
Replace-Text "130" "117"
Replace-Text "212" "182"

# Bold "2389" total manuscript count (standalone run, before "unique manuscripts.")
Replace-Text "2389" "2791"

# --- "Authors Reached" section ------------------------------------------
Replace-Text "4378 Authors Reached" "4973 Authors Reached"
Replace-Text "2389 manuscripts" "2791 manuscripts"
Replace-Text "1480 manuscripts have been published" "1724 manuscripts have been published"
Replace-Text "3154" "3631"
Replace-Text ", or about 2 authors per published article." ", or about 1.9 authors per published article."
Replace-Text "including the 614 articles not yet published" "including the 696 articles not yet published"
Replace-Text "4378" "4973"

# --- "Undergraduates" section --------------------------------------------
Replace-Text "166 Undergraduates" "201 Undergraduates"
Replace-Text "166 undergraduate students" "201 undergraduate students"

# --- Appendix: full roster of undergraduate students ----------------------
$oldRoster = "The following students have participated in the LDI Replication Lab, however briefly: Adam J. Faridi, Akshay Yadava, Albert Li, Alice Wei, Alizay Zartash, Amie Li, Ananya Bakshi, Andreas Psahos, Andres Aradillas Fernandez, Andrew Phiri, Andrew Wallace, Angelina Chen, Anjini Khanna, Anthony Peraza, Anurag Tiwari, Arnaav Sareen, Asha Patt, Ashley Cooray, Ashley Yu, Benjamin Zheng, Bianca Jimenez, Botao Yao, Cade Lenczycki, Caitlin Song, Christine Cho, Claire Vaughn, Connor Williams, Craig Schulman, Crystal Lim, Daniella Pena, David Wasser, Dhilan Bansal, Dmitry Shlyapnikov, Edward Vu, Eli Kolodezh, Elian Gomez, Elijah Ruiz, Elliott Serna, Emily Brydges, Emma Sbrollini, Ethan Carlson, Franklin Otieno, Gabriel Bond, Gary Wu, Gavin Bramley, Giles Mitchell, Hana Saleh, Hawi Tolera, Henry Tan, Hongyi Duan, Huey Le, Hyuk Son, Ilanith Nizard, Ilona Khimey, Jack Donnellan, Jacob Brogdon, Jacob Recht, Jade Yang, Jaeyoung Shim, Jai Kishore Kumar Chandnani, Janet Malzahn, Jared Martin, Jason Lan, Jeong Hyun Lee, Jessica Rizzo, Jill Crosby, John Park, Jonah Huang, Jonathan Temkin, Joshua Passell, Joshua Wallace, Julia Zimmerman, Kai Wang, Kareena Stowers, Kate Chanpong, Kate Hofer, Katheryn Ding, Kevin Bao, Kirin Eicher, Kirubeal Wondimu, Kristine Li, Kushal Kumar Reddy, Lauren Stubbs, Leah Kim, Lei Huang, Leslie Geng, Liam Cushen, Lilly Thomalla, LinchenZhang, Lincy Chen, Louis Liu, Luis Lopez Cabrera, Luke Trautwein, Lydia Reiner, Manas Gogula, Manvir Chahal, Marina Starikovsky, Mary-Jo Ajiduah, Mateo Pesa, Matt Wang, Matthew LaFontaine, Maximilian Yap, Melanie Brown, Melanie Chen, Meredith Welch, Micere Mugweru, Milena Zhu, Miranda Zhou, Mona Wu, Nahid Hassan, Naomi Li, Nathan Maidi, Nehedin Juarez, Nguyen Vo, Nicholas Esarte, Nicholas Swan, Nishat Tasneem Peuly, Olivia Liu, Ololade Omotoba, Peter Sanchez, Phalguni Miraj, Qianyi Liu, Raymond Wang, Rubal Mistry, Ryan Ali, Sam Evans, Satya Datla, Sean Miranda, Seong Hwan Kim, Sharan Banerjee, Sharon Liu, Sherry Li, Siddhi Malvankar, Siyang Elaine Yu, Sohit Gurung, Sourabh Velaga, Steve Yeh, Surita Basu, Suvd Khaliun, Sylverie Herbert, Syon Verma, Talia Boehm, Tarangana Thapa, Taren Daniels, Tommy Wang, Tony Ford Jr, Tyler Irving, Valerie Setiawan, Vansh Gupta, Vedha Mahesh, Victoria Liu, Vidya Balaji, Weilun Shi, Weiting Shen, Xiangru Li, Ximei Shen, Yanyun Chen, Yicheng Yang, Yuchang Tian, Zebang Xu, Zechariah Karsana, Kayla Yang, Jiayin Song, Nick Cichoski, Sharon Lin, Zhaojiahong Zhu."

$newRoster = "The following students have participated in the LDI Replication Lab, however briefly: Adam J. Faridi, Adis Bajraktarevic, Aishah Abdul-Hakeem, Akshay Yadava, Albert Li, Alice Wei, Alizay Zartash, Amie Li, Ananya Bakshi, Andreas Psahos, Andres Aradillas Fernandez, Andrew Phiri, Andrew Wallace, Angelina Chen, Anjini Khanna, Anthony Peraza, Anurag Tiwari, Arnaav Sareen, Aryan Shah, Asha Patt, Ashley Cooray, Ashley Yu, Austin Grattan, Benjamin Zheng, Bianca Jimenez, Botao Yao, Cade Lenczycki, Caitlin Song, Cameron Kang, Carissa Zheng, Christine Cho, Claire Vaughn, Connor Williams, Craig Schulman, Crystal Lim, Dane Clunk, Daniella Pena, David Graham, David Wasser, Dhilan Bansal, Dmitry Shlyapnikov, Edward Vu, Eirdeena Zakirah Mohd Nor Al-Azzad, Eli Kolodezh, Elian Gomez, Elijah Ruiz, Elizabeth Chen, Elliott Serna, Emily Brydges, Emma Sbrollini, Esther Park, Ethan Carlson, Farhan Mashrur, Franklin Otieno, Gabriel Bond, Gary Wu, Gavin Bramley, Gemie Sonamai, Giles Mitchell, Grace Fang, Grace Yichi Zhang, Hana Saleh, Hawi Tolera, Henry Tan, Hongyi Duan, Huey Le, Hyuk Son, Ilona Khimey, Jack Donnellan, Jacob Brogdon, Jacob Recht, Jade Yang, Jaeyoung Shim, Jai Kishore Kumar Chandnani, Janet Malzahn, Jared Martin, Jason Lan, Jeong Hyun Lee, Jessica Rizzo, Jesus Robles, Jiatian Yin, Jill Crosby, John Park, Jonah Huang, Jonathan Temkin, Joshua Passell, Joshua Wallace, Julia Zimmerman, Kai Wang, Kareena Stowers, Kate Chanpong, Kate Hofer, Katheryn Ding, Kevin Bao, Kirin Eicher, Kirubeal Wondimu, Kristine Li, Kushal Kumar Reddy, Lauren Stubbs, Leah Kim, Lei Huang, Leslie Geng, Liam Cushen, Lilly Thomalla, LinchenZhang, Lincy Chen, Lorraine Yang, Louis Liu, Luis Lopez Cabrera, Luke Trautwein, Lydia Reiner, Manas Gogula, Manvir Chahal, Marina Starikovsky, Mary-Jo Ajiduah, Mateo Pesa, Matt Wang, Matthew LaFontaine, Max Walder, Maximilian Yap, Melanie Brown, Melanie Chen, Meredith Welch, Micere Mugweru, Milena Zhu, Miranda Zhou, Mona Wu, Nahid Hassan, Naomi Li, Natalie Miller, Nathan Kwon, Nathan Maidi, Nehedin Juarez, Nguyen Vo, Nicholas Esarte, Nicholas Swan, Nishat Tasneem Peuly, Olivia Kim, Olivia Liu, Ololade Omotoba, Owen Kinney, Owen Ryan, Peter Sanchez, Phalguni Miraj, Pranav Jaganathan, Qianyi Liu, Ratchaphon Lertdamrongwong, Raymond Wang, Rishika Cherivirala, River Strumwasser, Rubal Mistry, Ryan Ali, Sam Evans, Sana Khan, Sanjana Bajaj, Satya Datla, Sean Miranda, Seong Hwan Kim, Sharan Banerjee, Sharon Liu, Sherry Li, Siddhi Malvankar, Siyang Elaine Yu, Sohit Gurung, Sonia Zhang, Sourabh Velaga, Steve Yeh, Surita Basu, Suvd Khaliun, Sylverie Herbert, Syon Verma, Talia Boehm, Talia Richmond, Tarangana Thapa, Taren Daniels, Tommy Wang, Tony Ford Jr, Tyler Irving, Valerie Setiawan, Vansh Gupta, Vedha Mahesh, Venna Huang, Victoria Liu, Vidya Balaji, Vivian Ye, Weilun Shi, Weiting Shen, William Howard, Wyatt Lewis, Xiangru Li, Ximei Shen, Yanyun Chen, Yicheng Yang, Yuchang Tian, Zebang Xu, Zechariah Karsana, Kayla Yang, Jiayin Song, Nick Cichoski, Sharon Lin, Zhaojiahong Zhu."

Replace-Text $oldRoster $newRoster

# --- Style cleanup ---------------------------------------------------------
# Title / TitleChar: drop the condensed-spacing + kerning tweak (equivalent to
# neutral/no-op values: no extra character spacing, no kerning threshold).
$titleStyle = $d.Styles.Item("Title")
$titleStyle.Font.Spacing = 0
$titleStyle.Font.Kerning = 0
$titleCharStyle = $d.Styles.Item("TitleChar")
$titleCharStyle.Font.Spacing = 0
$titleCharStyle.Font.Kerning = 0

# Author / Date: now based on Title (inheriting its centering), with an
# explicit, smaller 12pt run size.
$authorStyle = $d.Styles.Item("Author")
$authorStyle.BaseStyle = $titleStyle
$authorStyle.Font.Size = 12
$authorStyle.Font.SizeBi = 12

$dateStyle = $d.Styles.Item("Date")
$dateStyle.BaseStyle = $titleStyle
$dateStyle.Font.Size = 12
$dateStyle.Font.SizeBi = 12

Write-Output "done"
